$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.157093666666667
$ws.Range("H2").Value = 9.471281000000001
$ws.Range("I2").Value = 0.8981781966433163
$ws.Range("J2").Value = 0.8981781966433162
$ws.Range("M2").Value = 5.578493666666667
$ws.Range("N2").Value = 16.735481
$ws.Range("O2").Value = 0.1036332930693284
$ws.Range("P2").Value = 0.1036332930693284
$ws.Range("Q2").Value = 17.61182702457345
$ws.Range("R2").Value = 158.506443221161
$ws.Range("S2").Value = 0.09308116428121767
$ws.Range("T2").Value = 0.09308116428121768
$ws.Range("G3").Value = 3.157093666666667
$ws.Range("H3").Value = 9.471281000000001
$ws.Range("I3").Value = 0.8981781966433163
$ws.Range("J3").Value = 0.8981781966433162
$ws.Range("O3").Value = 0.06881911773528272
$ws.Range("P3").Value = 0.06881911773528274
$ws.Range("Q3").Value = 11.69537666555411
$ws.Range("R3").Value = 105.258389989987
$ws.Range("S3").Value = 0.06181183106206031
$ws.Range("T3").Value = 0.06181183106206031
$ws.Range("G4").Value = 3.157093666666667
$ws.Range("H4").Value = 9.471281000000001
$ws.Range("I4").Value = 0.8981781966433163
$ws.Range("J4").Value = 0.8981781966433162
$ws.Range("M4").Value = 24.77295966666667
$ws.Range("N4").Value = 74.31887900000001
$ws.Range("O4").Value = 0.4602144490493554
$ws.Range("P4").Value = 0.4602144490493556
$ws.Range("Q4").Value = 78.21055406822212
$ws.Range("R4").Value = 703.8949866139992
$ws.Range("S4").Value = 0.4133545839163474
$ws.Range("T4").Value = 0.4133545839163476
$ws.Range("G5").Value = 3.157093666666667
$ws.Range("H5").Value = 9.471281000000001
$ws.Range("I5").Value = 0.8981781966433163
$ws.Range("J5").Value = 0.8981781966433162
$ws.Range("M5").Value = 0.4291063333333334
$ws.Range("N5").Value = 1.287319
$ws.Range("O5").Value = 0.007971632676749163
$ws.Range("P5").Value = 0.007971632676749165
$ws.Range("Q5").Value = 1.354728887293222
$ws.Range("R5").Value = 12.192559985639
$ws.Range("S5").Value = 0.007159946661905496
$ws.Range("T5").Value = 0.007159946661905497
$ws.Range("G6").Value = 3.157093666666667
$ws.Range("H6").Value = 9.471281000000001
$ws.Range("I6").Value = 0.8981781966433163
$ws.Range("J6").Value = 0.8981781966433162
$ws.Range("M6").Value = 19.34413
$ws.Range("N6").Value = 58.03239
$ws.Range("O6").Value = 0.3593615074692841
$ws.Range("P6").Value = 0.3593615074692842
$ws.Range("Q6").Value = 61.07123031017667
$ws.Range("R6").Value = 549.6410727915901
$ws.Range("S6").Value = 0.3227706707217853
$ws.Range("T6").Value = 0.3227706707217853
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.3579033333333333
$ws.Range("H7").Value = 1.07371
$ws.Range("I7").Value = 0.1018218033566837
$ws.Range("J7").Value = 0.1018218033566837
$ws.Range("M7").Value = 5.578493666666667
$ws.Range("N7").Value = 16.735481
$ws.Range("O7").Value = 0.1036332930693284
$ws.Range("P7").Value = 0.1036332930693284
$ws.Range("Q7").Value = 1.996561478278889
$ws.Range("R7").Value = 17.96905330451
$ws.Range("S7").Value = 0.01055212878811073
$ws.Range("T7").Value = 0.01055212878811073
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.3579033333333333
$ws.Range("H8").Value = 1.07371
$ws.Range("I8").Value = 0.1018218033566837
$ws.Range("J8").Value = 0.1018218033566837
$ws.Range("O8").Value = 0.06881911773528272
$ws.Range("P8").Value = 0.06881911773528274
$ws.Range("Q8").Value = 1.325844189352222
$ws.Range("R8").Value = 11.93259770417
$ws.Range("S8").Value = 0.007007286673222424
$ws.Range("T8").Value = 0.007007286673222425
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.3579033333333333
$ws.Range("H9").Value = 1.07371
$ws.Range("I9").Value = 0.1018218033566837
$ws.Range("J9").Value = 0.1018218033566837
$ws.Range("M9").Value = 24.77295966666667
$ws.Range("N9").Value = 74.31887900000001
$ws.Range("O9").Value = 0.4602144490493554
$ws.Range("P9").Value = 0.4602144490493556
$ws.Range("Q9").Value = 8.866324841232222
$ws.Range("R9").Value = 79.79692357109001
$ws.Range("S9").Value = 0.04685986513300801
$ws.Range("T9").Value = 0.04685986513300804
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.3579033333333333
$ws.Range("H10").Value = 1.07371
$ws.Range("I10").Value = 0.1018218033566837
$ws.Range("J10").Value = 0.1018218033566837
$ws.Range("M10").Value = 0.4291063333333334
$ws.Range("N10").Value = 1.287319
$ws.Range("O10").Value = 0.007971632676749163
$ws.Range("P10").Value = 0.007971632676749165
$ws.Range("Q10").Value = 0.1535785870544444
$ws.Range("R10").Value = 1.38220728349
$ws.Range("S10").Value = 0.0008116860148436677
$ws.Range("T10").Value = 0.0008116860148436679
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.3579033333333333
$ws.Range("H11").Value = 1.07371
$ws.Range("I11").Value = 0.1018218033566837
$ws.Range("J11").Value = 0.1018218033566837
$ws.Range("M11").Value = 19.34413
$ws.Range("N11").Value = 58.03239
$ws.Range("O11").Value = 0.3593615074692841
$ws.Range("P11").Value = 0.3593615074692842
$ws.Range("Q11").Value = 6.923328607433333
$ws.Range("R11").Value = 62.30995746689999
$ws.Range("S11").Value = 0.03659083674749888
$ws.Range("T11").Value = 0.03659083674749889
